$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.824.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "'2.905.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'527.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "'144.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").Value = "'2.912.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").Value = "'6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "'3.412.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "'60.848.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "'22.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.03%  "
$ws.Range("D17").Value = "'2.898.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").Value = "'11.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'353.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.88%  "
$ws.Range("D22").Value = "'6.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'5.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "'64.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "'0.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").Value = "'0.180"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("D28").Value = "'0.984"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "'0.0₃0868"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.91%  "
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "'19.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").Value = "'153.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").Value = "'4.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'5.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.76%  "
$ws.Range("D37").Value = "'0.995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.02%  "
$ws.Range("D38").Value = "'1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("D39").Value = "'37.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.16%  "
$ws.Range("D41").Value = "'0.652"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").Value = "'3.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").Value = "'2.281.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.59%  "
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").Value = "'20.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.21%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'4.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.41%  "
$ws.Range("D48").Value = "'0.0237"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "'0.0915"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").Value = "'18.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.18%  "
